## Applies the "pushed the 10th commit" edit:
##  - Paragraph "My 7th task"  : the " task" that was split across a
##    proofErr-wrapped run collapses back into a single plain run.
##  - Paragraph "My 8th task"  : same collapse as above.
##  - Paragraph "My 9th task"  : the leading "My " + "9" runs merge into a
##    single "My 9" run, and its trailing " task" run splits apart again
##    into " " + proofErr(gramStart) + "task" + proofErr(gramEnd), i.e. it
##    now carries the split/flagged shape the two paragraphs above used to
##    have.
##  - A brand new paragraph "My 10th task" is appended, mirroring the shape
##    the old "My 9th task" paragraph had before this edit (two runs for
##    "My " / "10", superscript "th", merged " task").
##  - One extra empty paragraph is added at the end of the document.

$d = $word.ActiveDocument

# --- Step 1: rewrite paragraphs 7-9 and insert the new "My 10th task" ----
# Paragraph 6 is "My 6th task" (left untouched) and paragraph 9 is the old
# "My 9th task" paragraph. Replacing the whole span that covers paragraphs
# 7 through 9 in one go (rather than paragraph-by-paragraph) keeps the
# w:proofErr grammar markers from being left behind at the edges of the
# range, which is what happens if a paragraph containing one is edited in
# isolation.
$p6 = $d.Paragraphs.Item(6)
$p9 = $d.Paragraphs.Item(9)
$span = $d.Range($p6.Range.End, $p9.Range.End)

$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>My 7</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t xml:space="preserve"> task</w:t></w:r></w:p><w:p><w:r><w:t>My 8</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t xml:space="preserve"> task</w:t></w:r></w:p><w:p><w:r><w:t>My 9</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>task</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t xml:space="preserve">My </w:t></w:r><w:r><w:t>10</w:t></w:r><w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>th</w:t></w:r><w:r><w:t xml:space="preserve"> task</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$span.InsertXML($newParagraphsXml)

# --- Step 2: add one more empty paragraph at the very end of the body ----
# Spanning across the last two (already-empty) paragraphs and writing them
# back out as three empty paragraphs grows the trailing run of empty
# paragraphs by one, without disturbing the document's final paragraph
# mark (which cannot itself be replaced).
$secondLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$last = $d.Paragraphs.Last
$tailSpan = $d.Range($secondLast.Range.Start, $last.Range.End)

$extraEmptyParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$tailSpan.InsertXML($extraEmptyParaXml)
